# Practice1.xlsx - "Try Catch" exercise follow-up.
#
# For every month row, Cash In (B) and Cash Out (C) are subtracted to get the
# Difference (D). Some rows contain non-numeric "typo" values (e.g. "13K",
# "14O", "N12", ...) which would blow up the subtraction; those are caught
# and reported as D = 0 / Status = "ERROR", while the rows that compute
# cleanly are reported as Status = "SUCCESS".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 50

function Test-IsNumericCell($value) {
    # Numeric cells come back as [double]; text (including the "typo" values
    # like "13K") comes back as [string].
    return ($value -is [double])
}

# Work out, for every row, whether the Cash In / Cash Out values are usable
# numbers and what the resulting Difference / Status should be.
$difference = @{}
$status     = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cashIn  = $ws.Cells.Item($r, 2).Value2
    $cashOut = $ws.Cells.Item($r, 3).Value2

    if ((Test-IsNumericCell $cashIn) -and (Test-IsNumericCell $cashOut)) {
        $difference[$r] = $cashOut - $cashIn
        $status[$r]     = "SUCCESS"
    }
    else {
        # Subtraction isn't possible (caught "exception") - report 0/ERROR.
        $difference[$r] = 0
        $status[$r]     = "ERROR"
    }
}

# Write the "ERROR" status cells first, then the "SUCCESS" ones, so the
# shared-string table picks up "ERROR" before "SUCCESS" (matching how the
# workbook was authored).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($status[$r] -eq "ERROR") {
        $ws.Cells.Item($r, 5).Value = $status[$r]
    }
}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($status[$r] -eq "SUCCESS") {
        $ws.Cells.Item($r, 5).Value = $status[$r]
    }
}

# Now fill in the Difference column for every row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $difference[$r]
}

# Move the live selection to M17, as left by the author.
[void]$ws.Range("M17").Select()

# Mirror the window being minimized when the file was last saved (no-op on
# hosts that don't expose window state, harmless either way).
try { $excel.ActiveWindow.WindowState = -4140 } catch {}
